# Add a new "Italy" tab to the workbook, modeled on the existing
# "Slovakia" tab (same layout/template), then fill in the Italy-specific
# test data and leave Italy as the active/selected sheet - mirroring the
# manual edit captured in the commit "Test data added for Italy".

$wb = $excel.ActiveWorkbook

# The template sheet ("Slovakia") is the last sheet in the workbook.
$template = $wb.Worksheets.Item($wb.Worksheets.Count)

# Clone it and place the copy right after it; this becomes the new
# last sheet ("Slovakia (2)" for now - renamed below).
$template.Copy($null, $template)
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Fill in the Jira/ticket reference first, then the market name, so the
# shared-string table picks up the two new strings in that order.
$italy.Range("B4").ClearFormats()
$italy.Range("B4").Value = "NGC-3145/T2154"
$italy.Range("B2").Value = "Italy Market"

# Leave the cursor on B4 of the new sheet, matching the captured selection.
$italy.Range("B4").Select() | Out-Null

# The template sheet loses its "selected" state and reverts to a
# select-all-cells selection, as happens when focus moves off a tab.
$template.Cells.Select() | Out-Null

# Make Italy the active/visible tab.
$italy.Activate() | Out-Null
